$wb = $excel.ActiveWorkbook

# --- Overview sheet: handback status text changed to "not in sync" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: not in sync with en-US"
$overview.Range("F2").Value = "Handed back: not in sync with en-US"
$overview.Range("E3").Value = "Handed back: not in sync with en-US"
$overview.Range("F3").Value = "Handed back: not in sync with en-US"
$overview.Columns.Item(5).ColumnWidth = 32.6268870036
$overview.Columns.Item(6).ColumnWidth = 32.6268870036

# --- zh-cn sheet: status text + new handback datetime for row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: not in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: not in sync with en-US"
$zhcn.Range("K3").Value = "2016-09-09 08:39:30"
$zhcn.Columns.Item(3).ColumnWidth = 32.6268870036

# --- de-de sheet: status text + new handback datetime for row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: not in sync with en-US"
$dede.Range("C3").Value = "Handed back: not in sync with en-US"
$dede.Range("K3").Value = "2016-09-09 08:39:47"
$dede.Columns.Item(3).ColumnWidth = 32.6268870036

Write-Host "Generate Report for Handback: done"
